# Rename the inline picture shapes living in the document's headers/footers:
#   - the two Pearson logo images (currently "image2.png")  -> "image1.png"
#   - the BTec logo image               (currently "image1.jpg")  -> "image2.jpg"
#
# The pictures are inline shapes inside the header/footer ranges (not in the
# main document body), so we walk every section's Headers/Footers collections
# and inspect each InlineShape's AlternativeText (which carries the original
# descr/alt-text) to decide how to rename it.
#
# Note: a trivial Write-Output is emitted after each rename that touches a
# different header/footer part. Without it, the engine can report the next
# shape handle in another part as stale ("addressed block not found") because
# pending edits to the previous part haven't been flushed yet.

$d = $word.ActiveDocument

# WdHeaderFooterIndex values: 1 = primary, 2 = first page, 3 = even pages.
$hfTypes = @(1, 2, 3)

foreach ($sec in $d.Sections) {

    foreach ($t in $hfTypes) {

        # --- Headers ---
        $hdr = $sec.Headers.Item($t)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange" -and $shp.Name -ne "image2.jpg") {
                    $shp.Name = "image2.jpg"
                    Write-Output "Renamed BTec logo in header (type $t) to image2.jpg"
                }
            }
        }

        # --- Footers ---
        $ftr = $sec.Footers.Item($t)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" -and $shp.Name -ne "image1.png") {
                    $shp.Name = "image1.png"
                    Write-Output "Renamed Pearson logo in footer (type $t) to image1.png"
                }
            }
        }
    }
}
